# Incorporate Flora function (parent of Bush and Tree) into the design
# rationale bullet list, and remove the "Number of Fruit" method bullet
# that referenced GroundInterface.
#
# Concretely (per the target diff):
#   1. Append a new run containing "\" right after the existing
#      "Number of Fruit function in Ground Interface" text, in the same
#      bullet paragraph.
#   2. Add a brand-new bullet paragraph (same ListParagraph style /
#      numbering) right after it, containing "Flora Class added".

$d = $word.ActiveDocument

# Locate the paragraph that ends with the "Number of Fruit..." bullet -
# it is the last paragraph in the body.
$lastParaIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($lastParaIndex)

if ($targetPara.Range.Text -notlike "Number of Fruit function in Ground Interface*") {
    # Fall back to a direct search in case document structure differs.
    $searchRange = $d.Content
    $searchRange.Find.Execute("Number of Fruit function in Ground Interface",
                               $true, $false, $false, $false, $false,
                               $true, 1, $false, "", 0)
    $targetPara = $searchRange.Paragraphs.Item(1)
}

$tail = $targetPara.Range
$tail.Collapse(0)            # wdCollapseEnd -> move to just before the paragraph mark
$tail.InsertAfter("\")       # new run: "\"

# Force the newly inserted backslash to live in its own run (matching the
# authored formatting) by (re)stamping its character formatting.
$tail.Font.Name = "Times New Roman"
$tail.Font.NameAscii = "Times New Roman"
$tail.Font.NameBi = "Times New Roman"
$tail.Font.Size = 12

# Insert a new paragraph after the bullet; it inherits the ListParagraph
# style + numbering (numId 1, ilvl 0) from the paragraph it splits from.
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.InsertAfter("Flora Class added")
$newRange.Font.Name = "Times New Roman"
$newRange.Font.NameAscii = "Times New Roman"
$newRange.Font.NameBi = "Times New Roman"
$newRange.Font.Size = 12
